# Insert a new "Docentes responsáveis:" row pair into the LOQ4099 syllabus
# sheet, right after the "Objectives:" entry (row 11) and before
# "Programa resumido:" (the old row 12).
#
# Layout convention on this sheet: column A holds a bold label, and columns
# B/C hold the (duplicated) value text - sometimes the label and its value
# share one row, sometimes they are split across two rows (as is the case
# here: label on its own row, value on the next row).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push every row from 12 downward two rows further down, opening up a gap
# at rows 12-13 while preserving formatting/heights of all later rows.
$ws.Rows("12:13").Insert()

# New label row (A only).
$ws.Range("A12").Value = "Docentes responsáveis:"

# New value row (B and C both carry the same text, matching the existing
# "duplicate into column C" convention used throughout the sheet).
$ws.Range("B13").Value = "210064 - Eduardo Rezende Triboni"
$ws.Range("C13").Value = "210064 - Eduardo Rezende Triboni"

# The row Insert() operation stamps every column of the two new rows with
# blank-but-styled cells (carrying over column A/B/C formatting). Drop the
# cells that should stay genuinely empty so the row layout matches the
# label-only / value-only pattern used elsewhere on this sheet.
$ws.Range("B12").Clear()
$ws.Range("C12").Clear()
$ws.Range("A13").Clear()
